$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historical_data")

# Row 14 was "Export" (all zero historical trade values); it becomes "Import"
# and is populated with the historical import series (2000-2023, cols C:Z).
# Row 15 was "Import" (single 44.35 value in C15); it becomes "Export"
# and is zeroed out across the whole 2000-2023 range (cols C:Z).

$ws.Range("B14").Value = "Import"
$ws.Range("B15").Value = "Export"

$importValues = @(44.35, 48.38, 50.6, 50.97, 45.63, 49.15, 44.99, 46.28, 40.03, 44.96, 44.16, 45.73, 43.1, 42.14, 43.72, 46.38, 37.03, 37.76, 43.9, 38.14, 32.200000000000003, 42.79, 42.99, 51.25)

for ($i = 0; $i -lt $importValues.Length; $i++) {
    $ws.Cells.Item(14, 3 + $i).Value = $importValues[$i]
}

for ($col = 3; $col -le 26; $col++) {
    $ws.Cells.Item(15, $col).Value = 0
}
